# Edit script for sujil_mani_msc_project_proposal.docx
# Applies three surgical run-level changes (see commit message: "minor
# changes in grammer"):
#  1. Paragraph containing "new findings that are useful" -> drop the
#     "that" and " are " runs so the sentence reads "...new findings
#     useful for authorities...".
#  2. "Nominal Data" paragraph -> merge the " " and "sex" runs into a
#     single " sex" run (no visible text change).
#  3. "Discrete Data" paragraph -> merge the "," and " " runs into a
#     single ", " run (no visible text change).
#
# We use Range.InsertXML scoped to the whole paragraph (minus its
# trailing paragraph mark) and hand it the exact desired run markup so
# that untouched runs keep their original boundaries / formatting
# rather than being renormalised the way a plain Find/Replace would.

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

function Replace-ParagraphXml($doc, $para, [string]$bodyXml) {
    $full = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r = $doc.Range($para.Range.Start, $para.Range.End - 1)
    $r.InsertXML($full)
}

# --- 1. "that are" removal -------------------------------------------------
$p1 = Get-ParagraphByText $d "new findings that are useful"
if ($null -eq $p1) { throw "Could not locate 'new findings that are useful' paragraph" }
$p1Xml = '<w:p w14:paraId="5FA629C9" w14:textId="3307E9E2" w:rsidR="00A76704" w:rsidRDefault="005A247F" w:rsidP="00064F33"><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>T</w:t></w:r><w:r><w:t>his</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00A157F2"><w:t>iteration</w:t></w:r><w:r><w:t xml:space="preserve"> is a deep dive into </w:t></w:r><w:r w:rsidR="00CE4025"><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="00A157F2"><w:t>suicide dataset</w:t></w:r><w:r><w:t xml:space="preserve"> for knowing</w:t></w:r><w:r w:rsidR="00C936EB"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="001C16AD"><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="002D6C41"><w:t xml:space="preserve">reasons </w:t></w:r><w:r w:rsidR="001C16AD"><w:t xml:space="preserve">behind </w:t></w:r><w:r w:rsidR="00CF6CE6"><w:t>suicides around the world</w:t></w:r><w:r w:rsidR="002D6C41"><w:t>.</w:t></w:r><w:r w:rsidR="00842F73"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00536443"><w:t>Though,</w:t></w:r><w:r w:rsidR="002D6C41"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F120DE"><w:t>several</w:t></w:r><w:r w:rsidR="00111A7C"><w:t xml:space="preserve"> studies related to </w:t></w:r><w:r w:rsidR="00057E99"><w:t>this topic had</w:t></w:r><w:r w:rsidR="00111A7C"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="000C15C9"><w:t xml:space="preserve">been </w:t></w:r><w:r w:rsidR="00111A7C"><w:t>done in the past</w:t></w:r><w:r w:rsidR="00C936EB"><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="005857D3"><w:t>example: -</w:t></w:r><w:r w:rsidR="002E0FDC"><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rId8" w:history="1"><w:r w:rsidR="00C936EB" w:rsidRPr="00736CBB"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/></w:rPr><w:t>covid-19 time</w:t></w:r></w:hyperlink><w:r w:rsidR="00BF2495"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> s</w:t></w:r><w:r w:rsidR="00BF2495" w:rsidRPr="00BF2495"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t>uicide</w:t></w:r><w:r w:rsidR="006C3815" w:rsidRPr="00BF2495"><w:rPr><w:rStyle w:val="Hyperlink"/><w:color w:val="auto"/><w:u w:val="none"/></w:rPr><w:t xml:space="preserve"> research</w:t></w:r><w:r w:rsidR="00C936EB"><w:t>)</w:t></w:r><w:r w:rsidR="00536443"><w:t>.</w:t></w:r><w:r w:rsidR="00111A7C"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0062235D"><w:t>This</w:t></w:r><w:r w:rsidR="00111A7C"><w:t xml:space="preserve"> study </w:t></w:r><w:r w:rsidR="00057E99"><w:t xml:space="preserve">is going to </w:t></w:r><w:r w:rsidR="00842F73"><w:t>make</w:t></w:r><w:r w:rsidR="00057E99"><w:t xml:space="preserve"> new findings </w:t></w:r><w:r w:rsidR="00842F73"><w:t>useful</w:t></w:r><w:r w:rsidR="00057E99"><w:t xml:space="preserve"> for </w:t></w:r><w:r w:rsidR="00407497"><w:t>authorities</w:t></w:r><w:r w:rsidR="00057E99"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00407497"><w:t>to</w:t></w:r><w:r w:rsidR="00842F73"><w:t xml:space="preserve"> mak</w:t></w:r><w:r w:rsidR="00407497"><w:t>e</w:t></w:r><w:r w:rsidR="00842F73"><w:t xml:space="preserve"> policies that could reduce the mortality rate in </w:t></w:r><w:r w:rsidR="00BF09B4"><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="00842F73"><w:t>future</w:t></w:r><w:r w:rsidR="00057E99"><w:t xml:space="preserve">. I will be looking at different aspects of suicide and predicting how much more suicides are going to happen in the coming years in different </w:t></w:r><w:r w:rsidR="002F329D"><w:t>countries</w:t></w:r><w:r w:rsidR="00057E99"><w:t>.</w:t></w:r></w:p>'
Replace-ParagraphXml $d $p1 $p1Xml

# --- 2. "Nominal Data" run merge -------------------------------------------
$p2 = Get-ParagraphByText $d "Nominal Data"
if ($null -eq $p2) { throw "Could not locate 'Nominal Data' paragraph" }
$p2Xml = '<w:p w14:paraId="35971BD3" w14:textId="146F8667" w:rsidR="00233278" w:rsidRDefault="00F651AE" w:rsidP="00233278"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Nominal Data</w:t></w:r><w:r w:rsidR="00233278" w:rsidRPr="0012165A"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0012165A"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:-</w:t></w:r><w:r w:rsidR="00A20531"><w:t>c</w:t></w:r><w:r w:rsidR="00233278" w:rsidRPr="00233278"><w:t>ountry</w:t></w:r><w:r w:rsidR="00233278"><w:t>, Year,</w:t></w:r><w:r w:rsidR="00233278" w:rsidRPr="00233278"><w:t xml:space="preserve"> sex</w:t></w:r><w:r w:rsidR="00233278"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00233278" w:rsidRPr="00233278"><w:t>age</w:t></w:r><w:r w:rsidR="00233278"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00233278" w:rsidRPr="00233278"><w:t>country-year</w:t></w:r><w:r w:rsidR="00233278"><w:t>,</w:t></w:r><w:r w:rsidR="00A20531"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00233278"><w:t xml:space="preserve">generation </w:t></w:r><w:r w:rsidR="0038162F"><w:t>:- all these variables are nominal categorical</w:t></w:r></w:p>'
Replace-ParagraphXml $d $p2 $p2Xml

# --- 3. "Discrete Data" run merge ------------------------------------------
$p3 = Get-ParagraphByText $d "Discrete"
if ($null -eq $p3) { throw "Could not locate 'Discrete' paragraph" }
$p3Xml = '<w:p w14:paraId="4B08BE77" w14:textId="5772A7AD" w:rsidR="001C2CB9" w:rsidRDefault="005B10A2" w:rsidP="001C2CB9"><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Discrete</w:t></w:r><w:r w:rsidRPr="005B10A2"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Data :- </w:t></w:r><w:r w:rsidRPr="005B10A2"><w:t>suicides_no</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="005B10A2"><w:t>population</w:t></w:r></w:p>'
Replace-ParagraphXml $d $p3 $p3Xml

Write-Host "Done"
